# Remove the "Word version of this document" bullet item under
# "Additional resources". This paragraph contained a hyperlink to a
# now-obsolete Word-version link; a PDF version link replaces/joins the
# remaining items elsewhere (handled by other edits), so here we simply
# delete the whole list-item paragraph (bullet + hyperlink + paragraph mark).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Word version of this document*") {
        $p.Range.Delete()
        break
    }
}
